$wb = $excel.ActiveWorkbook
$wsMain = $wb.Worksheets.Item(2)   # "Главные"
$wsLinear = $wb.Worksheets.Item(3) # "Линейные"

$newTimestamp = "2025-11-25 03:03:21"

# ---- sheet2 (Главные) data updates ----
# Row 2
$wsMain.Range("C2").Value = 29
$wsMain.Range("D2").Value = 623
$wsMain.Range("E2").Value = 261
$wsMain.Range("F2").Value = 362
$wsMain.Range("G2").Value = 21.48
$wsMain.Range("H2").Value = 9
$wsMain.Range("I2").Value = 12.48
$wsMain.Range("J2").Value = 113
$wsMain.Range("K2").Value = 136
# Row 4
$wsMain.Range("C4").Value = 22
$wsMain.Range("D4").Value = 350
$wsMain.Range("E4").Value = 149
$wsMain.Range("F4").Value = 201
$wsMain.Range("G4").Value = 15.91
$wsMain.Range("H4").Value = 6.77
$wsMain.Range("I4").Value = 9.14
$wsMain.Range("J4").Value = 72
$wsMain.Range("K4").Value = 88
# Row 6
$wsMain.Range("C6").Value = 28
$wsMain.Range("D6").Value = 469
$wsMain.Range("E6").Value = 209
$wsMain.Range("F6").Value = 260
$wsMain.Range("G6").Value = 16.75
$wsMain.Range("H6").Value = 7.46
$wsMain.Range("I6").Value = 9.29
$wsMain.Range("J6").Value = 97
$wsMain.Range("K6").Value = 110
# Row 22
$wsMain.Range("C22").Value = 20
$wsMain.Range("D22").Value = 426
$wsMain.Range("E22").Value = 187
$wsMain.Range("F22").Value = 239
$wsMain.Range("G22").Value = 21.3
$wsMain.Range("H22").Value = 9.35
$wsMain.Range("I22").Value = 11.95
$wsMain.Range("J22").Value = 76
$wsMain.Range("K22").Value = 82
# Row 24
$wsMain.Range("C24").Value = 29
$wsMain.Range("D24").Value = 475
$wsMain.Range("E24").Value = 215
$wsMain.Range("F24").Value = 260
$wsMain.Range("G24").Value = 16.38
$wsMain.Range("H24").Value = 7.41
$wsMain.Range("I24").Value = 8.97
$wsMain.Range("J24").Value = 105
$wsMain.Range("K24").Value = 115

# ---- sheet3 (Линейные) data updates ----
# Row 6
$wsLinear.Range("C6").Value = 17
$wsLinear.Range("D6").Value = 297
$wsLinear.Range("E6").Value = 139
$wsLinear.Range("F6").Value = 158
$wsLinear.Range("G6").Value = 17.47
$wsLinear.Range("H6").Value = 8.18
$wsLinear.Range("I6").Value = 9.29
$wsLinear.Range("J6").Value = 62
$wsLinear.Range("K6").Value = 74
# Row 12
$wsLinear.Range("C12").Value = 22
$wsLinear.Range("D12").Value = 388
$wsLinear.Range("E12").Value = 185
$wsLinear.Range("F12").Value = 203
$wsLinear.Range("G12").Value = 17.64
$wsLinear.Range("H12").Value = 8.41
$wsLinear.Range("I12").Value = 9.23
$wsLinear.Range("J12").Value = 85
$wsLinear.Range("K12").Value = 94
# Row 16
$wsLinear.Range("C16").Value = 26
$wsLinear.Range("D16").Value = 457
$wsLinear.Range("E16").Value = 214
$wsLinear.Range("F16").Value = 243
$wsLinear.Range("G16").Value = 17.58
$wsLinear.Range("H16").Value = 8.23
$wsLinear.Range("I16").Value = 9.35
$wsLinear.Range("J16").Value = 97
$wsLinear.Range("K16").Value = 104
# Row 21
$wsLinear.Range("C21").Value = 30
$wsLinear.Range("D21").Value = 585
$wsLinear.Range("E21").Value = 246
$wsLinear.Range("F21").Value = 339
$wsLinear.Range("G21").Value = 19.5
$wsLinear.Range("H21").Value = 8.2
$wsLinear.Range("I21").Value = 11.3
$wsLinear.Range("J21").Value = 113
$wsLinear.Range("K21").Value = 137
# Row 24
$wsLinear.Range("C24").Value = 29
$wsLinear.Range("D24").Value = 516
$wsLinear.Range("E24").Value = 207
$wsLinear.Range("F24").Value = 309
$wsLinear.Range("G24").Value = 17.79
$wsLinear.Range("H24").Value = 7.14
$wsLinear.Range("I24").Value = 10.66
$wsLinear.Range("J24").Value = 91
$wsLinear.Range("K24").Value = 117
# Row 25
$wsLinear.Range("C25").Value = 8
$wsLinear.Range("D25").Value = 143
$wsLinear.Range("E25").Value = 84
$wsLinear.Range("F25").Value = 59
$wsLinear.Range("G25").Value = 17.88
$wsLinear.Range("H25").Value = 10.5
$wsLinear.Range("I25").Value = 7.38
$wsLinear.Range("J25").Value = 42
$wsLinear.Range("K25").Value = 27

# ---- sheet2 timestamp updates (rows 2-26) ----
for ($r = 2; $r -le 26; $r++) {
    $wsMain.Range("AA$r").Value = $newTimestamp
}

# ---- sheet3 timestamp updates (rows 2-26) ----
for ($r = 2; $r -le 26; $r++) {
    $wsLinear.Range("AA$r").Value = $newTimestamp
}
